$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2644.2942
$ws.Range("I19").Value = 2300
$ws.Range("J19").Value = 2950.3333
$ws.Range("K19").Value = 2300
$ws.Range("L19").Value = 2950.3333
$ws.Range("M19").Value = -2125
$ws.Range("N19").Value = -3300.3333
$ws.Range("H33").Value = 1754.4348
$ws.Range("I33").Value = 1562.7142
$ws.Range("K33").Value = 1562.7142
$ws.Range("M33").Value = -1333.7142
$ws.Range("H43").Value = 38699.75
$ws.Range("I43").Value = 38699.75
$ws.Range("K43").Value = 38699.75
$ws.Range("M43").Value = -38630.75
$ws.Range("H86").Value = 5370.9546
$ws.Range("I86").Value = 5093.125
$ws.Range("K86").Value = 5093.125
$ws.Range("M86").Value = -3970.125
$ws.Range("H88").Value = 7199.4
$ws.Range("J88").Value = 7499.25
$ws.Range("L88").Value = 7499.25
$ws.Range("N88").Value = -8311.25
$ws.Range("H89").Value = 5370.9546
$ws.Range("I89").Value = 5093.125
$ws.Range("K89").Value = 25465.625
$ws.Range("M89").Value = -19849.625
$ws.Range("H91").Value = 7199.4
$ws.Range("J91").Value = 7499.25
$ws.Range("L91").Value = 7499.25
$ws.Range("N91").Value = -10307.25
$ws.Range("H103").Value = 1274
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("H116").Value = 8209
$ws.Range("I116").Value = 7377.3335
$ws.Range("K116").Value = 7377.3335
$ws.Range("M116").Value = -3935.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2588
$ws.Range("I2").Value = 1609.6
$ws.Range("K2").Value = 1609.6
$ws.Range("M2").Value = -1496.6
$ws.Range("H32").Value = 2624.8928
$ws.Range("I32").Value = 1192.2885
$ws.Range("K32").Value = 1192.2885
$ws.Range("M32").Value = -905.2885000000001
$ws.Range("H45").Value = 6959.2324
$ws.Range("I45").Value = 9771.5
$ws.Range("J45").Value = 2658.1177
$ws.Range("K45").Value = 9771.5
$ws.Range("L45").Value = 2658.1177
$ws.Range("M45").Value = -9394.5
$ws.Range("N45").Value = -3412.1177
$ws.Range("H74").Value = 2586.0557
$ws.Range("I74").Value = 2535.8667
$ws.Range("K74").Value = 2535.8667
$ws.Range("M74").Value = -1661.8667
$ws.Range("H77").Value = 2586.0557
$ws.Range("I77").Value = 2535.8667
$ws.Range("K77").Value = 12679.3335
$ws.Range("M77").Value = -8311.333500000001
$ws.Range("H101").Value = 80000
$ws.Range("J101").Value = 80000
$ws.Range("L101").Value = 80000
$ws.Range("N101").Value = -86490
$ws.Range("H102").Value = 3492.5557
$ws.Range("I102").Value = 3492.5557
$ws.Range("K102").Value = 3492.5557
$ws.Range("M102").Value = -1870.5557
$ws.Range("H110").Value = 2802.0833
$ws.Range("I110").Value = 2140.375
$ws.Range("J110").Value = 4125.5
$ws.Range("K110").Value = 2140.375
$ws.Range("L110").Value = 4125.5
$ws.Range("M110").Value = -95.375
$ws.Range("N110").Value = -8215.5
$ws.Range("H116").Value = 2588
$ws.Range("I116").Value = 1609.6
$ws.Range("K116").Value = 1609.6
$ws.Range("M116").Value = 684.4000000000001
$ws.Range("H122").Value = 4196.0713
$ws.Range("I122").Value = 4160.091
$ws.Range("K122").Value = 12480.273
$ws.Range("M122").Value = -10030.273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3950.8333
$ws.Range("I99").Value = 3810
$ws.Range("K99").Value = 3810
$ws.Range("M99").Value = -2312
$ws.Range("H105").Value = 3572.1428
$ws.Range("I105").Value = 3668
$ws.Range("J105").Value = 3332.5
$ws.Range("K105").Value = 3668
$ws.Range("L105").Value = 3332.5
$ws.Range("M105").Value = -1921
$ws.Range("N105").Value = -6826.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5070.15
$ws.Range("I31").Value = 5835.2144
$ws.Range("K31").Value = 5835.2144
$ws.Range("M31").Value = -5540.2144
$ws.Range("H34").Value = 5070.15
$ws.Range("I34").Value = 5835.2144
$ws.Range("K34").Value = 5835.2144
$ws.Range("M34").Value = -5633.2144
$ws.Range("H58").Value = 4023.6875
$ws.Range("I58").Value = 4127.1924
$ws.Range("K58").Value = 4127.1924
$ws.Range("M58").Value = -3924.1924
$ws.Range("H105").Value = 1403.4117
$ws.Range("I105").Value = 1426.8
$ws.Range("K105").Value = 1426.8
$ws.Range("M105").Value = 320.2
$ws.Range("H107").Value = 1134.909
$ws.Range("I107").Value = 998.5
$ws.Range("K107").Value = 998.5
$ws.Range("M107").Value = 921.5
$ws.Range("H122").Value = 123395.6
$ws.Range("I122").Value = 139604.1
$ws.Range("K122").Value = 418812.3
$ws.Range("M122").Value = -416362.3
$ws.Range("H134").Value = 6055.4517
$ws.Range("I134").Value = 5822.72
$ws.Range("J134").Value = 7025.1665
$ws.Range("K134").Value = 17468.16
$ws.Range("L134").Value = 21075.4995
$ws.Range("M134").Value = -14933.16
$ws.Range("N134").Value = -26145.4995
$ws.Range("H136").Value = 4023.6875
$ws.Range("I136").Value = 4127.1924
$ws.Range("K136").Value = 12381.5772
$ws.Range("M136").Value = -9831.5772
$ws.Range("H141").Value = 74427.14
$ws.Range("J141").Value = 71832.164
$ws.Range("L141").Value = 71832.164
$ws.Range("N141").Value = -82192.164

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 84282536
$ws.Range("I4").Value = 91944580
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 275833740
$ws.Range("L4").Value = 450
$ws.Range("M4").Value = -275833628
$ws.Range("N4").Value = -674
$ws.Range("H34").Value = 2312.24
$ws.Range("I34").Value = 142.16667
$ws.Range("J34").Value = 2997.5264
$ws.Range("K34").Value = 426.50001
$ws.Range("L34").Value = 8992.5792
$ws.Range("M34").Value = -342.50001
$ws.Range("N34").Value = -9160.5792
$ws.Range("H100").Value = 11249
$ws.Range("J100").Value = 11249
$ws.Range("L100").Value = 33747
$ws.Range("N100").Value = -35369

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 33000
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H70").Value = 6380.75
$ws.Range("I70").Value = 6435.143
$ws.Range("K70").Value = 6435.143
$ws.Range("M70").Value = -6165.143
$ws.Range("H73").Value = 6380.75
$ws.Range("I73").Value = 6435.143
$ws.Range("K73").Value = 6435.143
$ws.Range("M73").Value = -5499.143
$ws.Range("H97").Value = 585.6429000000001
$ws.Range("I97").Value = 584.53845
$ws.Range("K97").Value = 584.53845
$ws.Range("M97").Value = -88.53845000000001
$ws.Range("H106").Value = 49999
$ws.Range("J106").Value = 49999
$ws.Range("L106").Value = 49999
$ws.Range("N106").Value = -52523
$ws.Range("H132").Value = 5848.1875
$ws.Range("I132").Value = 5731.5454
$ws.Range("J132").Value = 6104.8
$ws.Range("K132").Value = 17194.6362
$ws.Range("L132").Value = 18314.4
$ws.Range("M132").Value = -14664.6362
$ws.Range("N132").Value = -23374.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15242.8
$ws.Range("I7").Value = 16307
$ws.Range("K7").Value = 16307
$ws.Range("M7").Value = -16195
$ws.Range("H82").Value = 1423.9
$ws.Range("I82").Value = 1087.5
$ws.Range("K82").Value = 1087.5
$ws.Range("M82").Value = -726.5
$ws.Range("H85").Value = 1423.9
$ws.Range("I85").Value = 1087.5
$ws.Range("K85").Value = 1087.5
$ws.Range("M85").Value = 160.5
$ws.Range("H126").Value = 15242.8
$ws.Range("I126").Value = 16307
$ws.Range("K126").Value = 48921
$ws.Range("M126").Value = -46451
$ws.Range("H132").Value = 6271.5454
$ws.Range("I132").Value = 5398.4
$ws.Range("J132").Value = 6999.1665
$ws.Range("K132").Value = 16195.2
$ws.Range("L132").Value = 20997.4995
$ws.Range("M132").Value = -13665.2
$ws.Range("N132").Value = -26057.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()
$ws.Range("H41").Value = 11999
$ws.Range("J41").Value = 11999
$ws.Range("L41").Value = 11999
$ws.Range("N41").Value = -12779
$ws.Range("H81").Value = 5141.778
$ws.Range("I81").Value = 4713.1665
$ws.Range("J81").Value = 5999
$ws.Range("K81").Value = 9426.333000000001
$ws.Range("L81").Value = 11998
$ws.Range("M81").Value = -8365.333000000001
$ws.Range("N81").Value = -14120
$ws.Range("H84").Value = 5141.778
$ws.Range("I84").Value = 4713.1665
$ws.Range("J84").Value = 5999
$ws.Range("K84").Value = 47131.665
$ws.Range("L84").Value = 59990
$ws.Range("M84").Value = -41827.665
$ws.Range("N84").Value = -70598
$ws.Range("H107").Value = 2659.0417
$ws.Range("I107").Value = 2585.8333
$ws.Range("K107").Value = 7757.499899999999
$ws.Range("M107").Value = -5837.499899999999
$ws.Range("H122").Value = 2724.3572
$ws.Range("I122").Value = 2750.4546
$ws.Range("K122").Value = 8251.363799999999
$ws.Range("M122").Value = -5801.363799999999
$ws.Range("H136").Value = 2598.6667
$ws.Range("I136").Value = 2339.6904
$ws.Range("J136").Value = 3323.8
$ws.Range("K136").Value = 7019.0712
$ws.Range("L136").Value = 9971.400000000001
$ws.Range("M136").Value = -4469.0712
$ws.Range("N136").Value = -15071.4
